$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (item 15) is missing data in columns G:J compared to the other
# fully-filled rows in the table (e.g. rows 12, 13, 20, 29, 30, 31).
# Copy the formatting of an already complete row so G18/H18 keep their
# existing style and I18/J18 pick up the correct "green fill + thick
# right border" style used throughout the sheet, then fill in the values.
$ws.Range("G12:J12").Copy() | Out-Null
$ws.Range("G18:J18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 5

# Update the active selection to the last edited cell, as recorded in the
# workbook view state.
$ws.Range("J18").Select() | Out-Null
